$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws4 = $wb.Worksheets.Item("Sheet4")

# ---------------------------------------------------------------------------
# "added new quantities into excel file sheet1 and 2" / "fixed typo for Q1
# and Q2" -- reallocate quantities on Sheet2 rows 5, 8 and 12 (row totals are
# preserved, only the split across destination columns changes).

# Row 5: move 1,000,000 out of C and 4,000,000 out of G into the new H
# quantity (5,000,000 shipped instead of 0).
$ws2.Range("C5").Value = 0
$ws2.Range("G5").Value = 7000000
$ws2.Range("H5").Value = 5000000

# Row 8: shift 1,000,000 from D into C.
$ws2.Range("C8").Value = 4000000
$ws2.Range("D8").Value = 1000000

# Row 12: add a new 1,000,000 quantity in D, move the 5,000,000 previously in
# H into G instead.
$ws2.Range("D12").Value = 1000000
$ws2.Range("G12").Value = 17000000
$ws2.Range("H12").Value = 0

# ---------------------------------------------------------------------------
# Restore/update the saved view state (active sheet + per-sheet selection)
# to match what was recorded when the workbook was last saved.

$ws1.Select()
$ws1.Range("D15").Select()
# Sheet1 was scrolled so row 13 is the top visible row.
try { $excel.ActiveWindow.ScrollRow = 13 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}

$ws4.Select()
$ws4.Range("D15").Select()

# Sheet3's selection is left untouched (A1:H13 / active cell H13), matching
# the unchanged state in the workbook.

# Sheet2 becomes the active/visible tab when the workbook is saved.
$ws2.Select()
$ws2.Range("H13").Select()
